$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.467.14'
$ws.Range("E2").Value = '  -2.10%  '
$ws.Range("D3").Value = '3.158.91'
$ws.Range("E3").Value = '  -6.18%  '
$ws.Range("D5").Value = '''205.43'
$ws.Range("E5").Value = '  -6.75%  '
$ws.Range("D6").Value = '''608.31'
$ws.Range("E6").Value = '  -6.33%  '
$ws.Range("E7").Value = '  -9.11%  '
$ws.Range("D8").Value = '''0.660'
$ws.Range("E8").Value = '  +7.25%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '3.150.44'
$ws.Range("E10").Value = '  -6.21%  '
$ws.Range("D11").Value = '''0.534'
$ws.Range("E11").Value = '  -14.05%  '
$ws.Range("E12").Value = '  +5.47%  '
$ws.Range("E13").Value = '  -17.01%  '
$ws.Range("D14").Value = '3.736.48'
$ws.Range("E14").Value = '  -6.16%  '
$ws.Range("D15").Value = '''5.24'
$ws.Range("E15").Value = '  -5.29%  '
$ws.Range("D16").Value = '87.084.40'
$ws.Range("E16").Value = '  -2.49%  '
$ws.Range("D17").Value = '''32.02'
$ws.Range("E17").Value = '  -11.90%  '
$ws.Range("D18").Value = '3.148.21'
$ws.Range("E18").Value = '  -6.40%  '
$ws.Range("E19").Value = '  -4.48%  '
$ws.Range("D20").Value = '''13.37'
$ws.Range("E20").Value = '  -9.82%  '
$ws.Range("D21").Value = '''413.66'
$ws.Range("E21").Value = '  -9.78%  '
$ws.Range("D22").Value = '''8.47'
$ws.Range("E22").Value = '  -12.82%  '
$ws.Range("D23").Value = '''5.08'
$ws.Range("E23").Value = '  -8.80%  '
$ws.Range("D24").Value = '''5.16'
$ws.Range("E24").Value = '  -6.98%  '
$ws.Range("D25").Value = '''11.80'
$ws.Range("E25").Value = '  -7.82%  '
$ws.Range("D26").Value = '3.316.18'
$ws.Range("E26").Value = '  -5.91%  '
$ws.Range("D27").Value = '''73.28'
$ws.Range("E27").Value = '  -7.65%  '
$ws.Range("D28").Value = '''0.0000130'
$ws.Range("E28").Value = '  -9.98%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -19.01%  '
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").Value = '''537.10'
$ws.Range("E32").Value = '  -10.07%  '
$ws.Range("D33").Value = '''8.26'
$ws.Range("E33").Value = '  -12.21%  '
$ws.Range("D34").Value = '''1.31'
$ws.Range("E34").Value = '  -17.40%  '
$ws.Range("D35").Value = '''6.69'
$ws.Range("E35").Value = '  -8.68%  '
$ws.Range("E36").Value = '  -13.02%  '
$ws.Range("D37").Value = '''0.133'
$ws.Range("E37").Value = '  -8.36%  '
$ws.Range("D38").Value = '''21.87'
$ws.Range("E38").Value = '  -7.22%  '
$ws.Range("D39").Value = '''21.80'
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("D40").Value = '''0.997'
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").Value = '''2.99'
$ws.Range("E41").Value = '  -6.53%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = '''1.89'
$ws.Range("E43").Value = '  -12.35%  '
$ws.Range("D44").Value = '''0.371'
$ws.Range("E44").Value = '  -13.84%  '
$ws.Range("D45").Value = '''147.48'
$ws.Range("E45").Value = '  -6.55%  '
$ws.Range("D46").Value = '''172.00'
$ws.Range("E46").Value = '  -8.77%  '
$ws.Range("D47").Value = '''43.20'
$ws.Range("E47").Value = '  -6.99%  '
$ws.Range("E48").Value = '  +4.82%  '
$ws.Range("D49").Value = '''1.25'
$ws.Range("E49").Value = '  -14.68%  '
$ws.Range("D50").Value = '''3.96'
$ws.Range("E50").Value = '  -12.15%  '
$ws.Range("D51").Value = '''0.693'
$ws.Range("E51").Value = '  -11.97%  '
